# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45177
